# DLP woordenlijst.xlsx — "Wat vertalingen erdoorheen gejast"
# Adds a handful of new English/Dutch term pairs to the glossary table,
# adds a new explanatory "Kolom1" column, tidies up a stray empty styled
# cell, and adds a clarifying remark about the "Raise Levies" entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Row 80: drop the stray empty, fill-styled cell in column C -----------
$ws.Range("C80").Clear()

# --- New glossary entries appended as rows 124-127 -------------------------
$ws.Range("A124").Value = "Alliance"
$ws.Range("B124").Value = "Bondgenootschap"
$ws.Range("D124").Value = "2024-05-28"

$ws.Range("A125").Value = "Allies"
$ws.Range("B125").Value = "Bondgenoten"
$ws.Range("D125").Value = "2024-05-28"

$ws.Range("A126").Value = "Allied"
$ws.Range("B126").Value = "Verbonden"
$ws.Range("D126").Value = "2024-05-28"

$ws.Range("A127").Value = "Concubine"
$ws.Range("B127").Value = "Bijvrouw"
$ws.Range("D127").Value = "2024-05-28"

# --- Row 92 ("Raise Levies" / "Baanderheren Oproepen"): add a remark ------
$ws.Range("F92").Value = "Heb zelf de context hiervan ietwat aangepast. Gezien het de Baanderheren waren die moesten worden opgeroepen en zij kwamen dan met de Heervaarders"

# --- Extend the table with a new last column ("Kolom1") -------------------
$newCol = $lo.ListColumns.Add()
$ws.Range("F1").Value = "Kolom1"

# --- Reflect the author's final cursor position ----------------------------
$ws.Range("F93").Select() | Out-Null
